$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: lat -> latitude, lng -> longitude
$ws.Range("G1").Value = "latitude"
$ws.Range("H1").Value = "longitude"

# Move the selection to I1 (as reflected in the saved sheet view)
$ws.Range("I1").Select()
